$d = $word.ActiveDocument

$pairs = @(
    @("466×8=", "542×7="),
    @("609×3=", "695×4="),
    @("411×5=", "681×3="),
    @("550×8=", "341×7="),
    @("959×9=", "229×3="),
    @("609×6=", "742×9="),
    @("945×3=", "826×4="),
    @("907×4=", "739×3="),
    @("143×8=", "464×4="),
    @("838×3=", "878×2="),
    @("563×5=", "158×7="),
    @("195×6=", "506×2="),
    @("943×7=", "976×5="),
    @("383×3=", "401×6="),
    @("524×5=", "399×4="),
    @("558×7=", "375×3="),
    @("396×4=", "822×7="),
    @("396×5=", "253×7="),
    @("273×8=", "914×6="),
    @("660×2=", "362×5="),
    @("241×5=", "342×4="),
    @("424×2=", "428×4="),
    @("136×6=", "217×8="),
    @("970×6=", "621×5="),
    @("483×3=", "538×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
